$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2: status moved back to "In Progress"; update-date refreshed ---
$ws.Range("D2").Value = "In Progress"
$ws.Range("F2").Value = "2023-08-02 16:31:58"

# --- Row 4: status closed; update-date refreshed ---
$ws.Range("D4").Value = "Closed"
$ws.Range("F4").Value = "2023-08-02 16:57:38"

# --- Row 7: update note edited; status closed; update-date refreshed ---
$ws.Range("C7").Value = "fgfg clossed"
$ws.Range("D7").Value = "Closed"
$ws.Range("F7").Value = "2023-08-02 16:56:51"

# --- Row 8: an "Update" note was posted via the popup, reopening the bug ---
$ws.Range("C8").Value = "printer is working"
$ws.Range("D8").Value = "In Progress"
$ws.Range("F8").Value = "2023-08-02 16:58:24"

# --- New bug reports submitted through the view-bug "update" popup ---
# Row 9
$ws.Range("A9").Value = 1
$ws.Range("B9").Value = "fan not working"
$ws.Range("C9").Value = "fan not working"
$ws.Range("D9").Value = "In Progress"
$ws.Range("E9").Value = "2023-08-02 17:09:57"
$ws.Range("F9").Value = "2023-08-02 17:10:14"
$ws.Range("A8").Copy()
$ws.Range("A9").PasteSpecial(-4122)

# Row 10
$ws.Range("A10").Value = 3
$ws.Range("B10").Value = "printer is not working"
$ws.Range("C10").Value = "printer is not working"
$ws.Range("D10").Value = "In Progress"
$ws.Range("E10").Value = "2023-08-02 17:10:04"
$ws.Range("F10").Value = "2023-08-02 17:10:27"
$ws.Range("A8").Copy()
$ws.Range("A10").PasteSpecial(-4122)

$excel.CutCopyMode = $false
